# feat: add 2022-Q4 data
#
# The workbook tracks this stock's holders by quarter. A new "2022-Q4" sheet
# is introduced (placed right after the "总计" summary sheet, before the
# existing "2022-Q3" sheet), and the "总计" summary sheet gets a new row for
# the 2022-Q4 totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q4" worksheet by duplicating "2022-Q3" (so it
# inherits the same column layout/styles), inserted immediately before it.
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3)
$wsQ4 = $wb.Worksheets.Item("2022-Q3 (2)")
$wsQ4.Name = "2022-Q4"

# The Q3 snapshot had 3 funds (000041, 519601, 519602); the Q4 snapshot only
# has the two 海富通 funds, so drop the 000041 row before filling in values.
$wsQ4.Rows(2).Delete()

# Row 2 -> fund 519601
$wsQ4.Range("A2").Value = 0
$wsQ4.Range("B2:G2").NumberFormat = "@"
$wsQ4.Range("B2").Value = "519601"
$wsQ4.Range("C2").Value = "海富通中国海外精选混合（QDII）"
$wsQ4.Range("D2").Value = "0.58"
$wsQ4.Range("E2").Value = "94.54"
$wsQ4.Range("F2").Value = "3.90"
$wsQ4.Range("G2").Value = "0.0226"
$wsQ4.Range("H2").Value = 6

# Row 3 -> fund 519602
$wsQ4.Range("A3").Value = 1
$wsQ4.Range("B3:G3").NumberFormat = "@"
$wsQ4.Range("B3").Value = "519602"
$wsQ4.Range("C3").Value = "海富通大中华精选混合（QDII）"
$wsQ4.Range("D3").Value = "0.11"
$wsQ4.Range("E3").Value = "89.56"
$wsQ4.Range("F3").Value = "3.69"
$wsQ4.Range("G3").Value = "0.0041"
$wsQ4.Range("H3").Value = 6

# ---------------------------------------------------------------------------
# Step 2: add the 2022-Q4 row to the "总计" (summary) sheet, pushing the
# existing 2022-Q3 / 2022-Q2 rows down by one.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# extend column-A's styling down onto the new row 4 (matches rows 2 & 3)
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)

# row 4: 2022-Q2 (unchanged data, shifted down from row 3)
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q2"
$wsTotal.Range("C4").Value = 3
$wsTotal.Range("D4").Value = 0.44

# row 3: 2022-Q3 (unchanged data, shifted down from row 2)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 3
$wsTotal.Range("D3").Value = 0.42

# row 2: new 2022-Q4 entry
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.03

# restore the original active tab (2022-Q2, the last sheet) -- copying a
# sheet shifts the selection onto the newly created copy otherwise.
$wb.Worksheets.Item("2022-Q2").Activate()

Write-Host "2022-Q4 sheet added"
